$wb = $excel.ActiveWorkbook

# "Metadata" sheet holds the value set properties
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/rx-mail-or-retail"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"
$meta.Range("B11").Value = "LinuxForHealth standard value set indicating the purchase place of the prescription."

# "Include from Rx Mail Or Retai" sheet holds the code system reference
$codes = $wb.Worksheets.Item("Include from Rx Mail Or Retai")

$codes.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/rx-mail-or-retail"
